$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "40"
$ws.Range("C9").Characters(27, 9).Text = "10/2/2023"
$ws.Range("C9").Characters(47, 9).Text = "10/8/2023"

# --- Weekly crime statistics table updates (rows 16-29) ---
# Row 16
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = 2.04081632653
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -35.897435897435
$ws.Range("N16").Value = -82.394366197183

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -44.444444444444
$ws.Range("J17").Value = 76
$ws.Range("K17").Value = 9.210526315789
$ws.Range("L17").Value = 29.6875
$ws.Range("M17").Value = 9.210526315789
$ws.Range("N17").Value = -65.560165975103

# Row 18
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 300
$ws.Range("L18").Value = 32.142857142857
$ws.Range("M18").Value = -12.941176470588
$ws.Range("N18").Value = -77.439024390243

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0"
$ws.Range("C17").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "***.*"
$ws.Range("C17").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("F19").Value = 12
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -29.411764705882
$ws.Range("I19").Value = 125
$ws.Range("K19").Value = -20.886075949367
$ws.Range("L19").Value = 31.578947368421
$ws.Range("M19").Value = -28.977272727272
$ws.Range("N19").Value = -31.693989071038

# Row 20
$ws.Range("I15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 50
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = 11.111111111111
$ws.Range("L20").Value = 25
$ws.Range("M20").Value = -1.960784313725
$ws.Range("N20").Value = -83.974358974359

# Row 21
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 60
$ws.Range("F21").Value = 31
$ws.Range("G21").Value = 40
$ws.Range("H21").Value = -22.5
$ws.Range("I21").Value = 388
$ws.Range("J21").Value = 418
$ws.Range("K21").Value = -7.177033492822
$ws.Range("L21").Value = 35.191637630662
$ws.Range("M21").Value = -17.27078891258
$ws.Range("N21").Value = -71.533382245047

# Row 22
$ws.Range("I15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("K15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = -66.666666666666

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -63.636363636363
$ws.Range("J23").Value = 93
$ws.Range("K23").Value = -8.602150537634

# Row 24
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = 57.5
$ws.Range("I24").Value = 476
$ws.Range("J24").Value = 347
$ws.Range("K24").Value = 37.175792507204
$ws.Range("L24").Value = 70.609318996415
$ws.Range("M24").Value = 37.175792507204

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 30.76923076923
$ws.Range("I25").Value = 127
$ws.Range("J25").Value = 143
$ws.Range("K25").Value = -11.188811188811
$ws.Range("L25").Value = 1.6
$ws.Range("M25").Value = -40.930232558139

# Row 26
$ws.Range("I15").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 1
$ws.Range("I26").Value = 9
$ws.Range("K26").Value = 28.571428571428
$ws.Range("L26").Value = 50

# Row 27
$ws.Range("L27").Value = -47.368421052631

# Row 28
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
$ws.Range("N28").Value = -85.714285714285

# Row 29
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = -82.051282051282

$excel.CutCopyMode = 0